$p = $ppt.ActivePresentation

# Slide 4: "Describing the business logic" -> "1. Describing the business logic"
$p.Slides.Item(4).Shapes.Item(1).TextFrame.TextRange.Text = "1. Describing the business logic"

# Slide 5: "Describing the business logic" -> "1. Describing the business logic"
$p.Slides.Item(5).Shapes.Item(1).TextFrame.TextRange.Text = "1. Describing the business logic"

# Slide 6: "Describing business logic" -> "1. Describing business logic"
$p.Slides.Item(6).Shapes.Item(1).TextFrame.TextRange.Text = "1. Describing business logic"

# Slide 7: "Adapter layer" -> "2. Adapter layer"
$p.Slides.Item(7).Shapes.Item(1).TextFrame.TextRange.Text = "2. Adapter layer"

# Slide 8: "Set up adapter layer " -> "2. Set up adapter layer "
$p.Slides.Item(8).Shapes.Item(1).TextFrame.TextRange.Text = "2. Set up adapter layer "

# Slide 9: "Data Binding " -> "3. Data Binding "
$p.Slides.Item(9).Shapes.Item(1).TextFrame.TextRange.Text = "3. Data Binding "

# Slide 10: "Binding data" -> "3. Binding data"
$p.Slides.Item(10).Shapes.Item(1).TextFrame.TextRange.Text = "3. Binding data"

# Slide 11: "Execute!" -> "4. Execute!"
$p.Slides.Item(11).Shapes.Item(1).TextFrame.TextRange.Text = "4. Execute!"

# Slide 12: "Execute!" -> "4. Execute!"
$p.Slides.Item(12).Shapes.Item(1).TextFrame.TextRange.Text = "4. Execute!"
